$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 26663
$ws.Cells.Item(2, 5).Value = 1268
$ws.Cells.Item(2, 6).Value = 1268
$ws.Cells.Item(2, 7).Value = 1401
$ws.Cells.Item(2, 8).Value = 1020
$ws.Cells.Item(2, 9).Value = 1016
$ws.Cells.Item(2, 10).Value = 4
$ws.Cells.Item(2, 11).Value = 18439
$ws.Cells.Item(2, 12).Value = 9612
$ws.Cells.Item(2, 13).Value = 8827
$ws.Cells.Item(2, 14).Value = 8785
$ws.Cells.Item(2, 15).Value = 42
$ws.Cells.Item(2, 16).Value = 230
$ws.Cells.Item(2, 17).Value = 559
$ws.Cells.Item(2, 18).Value = -1273
$ws.Cells.Item(2, 19).Value = 1299
$ws.Cells.Item(2, 20).Value = 390
$ws.Cells.Item(2, 21).Value = 169
$ws.Cells.Item(2, 22).Value = 184
$ws.Cells.Item(2, 23).Value = 4.75
$ws.Cells.Item(2, 24).Value = 3.83
$ws.Cells.Item(2, 25).Value = 13.19
$ws.Cells.Item(2, 26).Value = 5.75
$ws.Cells.Item(2, 27).Value = 108.89
$ws.Cells.Item(2, 28).Value = 4237.27
$ws.Cells.Item(2, 29).Value = 883
$ws.Cells.Item(2, 30).Value = 19.47
$ws.Cells.Item(2, 31).Value = 8121
$ws.Cells.Item(2, 32).Value = 2.12
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 115041225

# Row 3
$ws.Cells.Item(3, 4).Value = 28067
$ws.Cells.Item(3, 5).Value = 1272
$ws.Cells.Item(3, 6).Value = 1272
$ws.Cells.Item(3, 7).Value = 1495
$ws.Cells.Item(3, 8).Value = 817
$ws.Cells.Item(3, 9).Value = 795
$ws.Cells.Item(3, 10).Value = 22
$ws.Cells.Item(3, 11).Value = 18669
$ws.Cells.Item(3, 12).Value = 10665
$ws.Cells.Item(3, 13).Value = 8003
$ws.Cells.Item(3, 14).Value = 7859
$ws.Cells.Item(3, 15).Value = 144
$ws.Cells.Item(3, 16).Value = 230
$ws.Cells.Item(3, 17).Value = 1024
$ws.Cells.Item(3, 18).Value = 988
$ws.Cells.Item(3, 19).Value = -1628
$ws.Cells.Item(3, 20).Value = 177
$ws.Cells.Item(3, 21).Value = 847
$ws.Cells.Item(3, 22).Value = 270
$ws.Cells.Item(3, 23).Value = 4.53
$ws.Cells.Item(3, 24).Value = 2.91
$ws.Cells.Item(3, 25).Value = 9.550000000000001
$ws.Cells.Item(3, 26).Value = 4.41
$ws.Cells.Item(3, 27).Value = 133.26
$ws.Cells.Item(3, 28).Value = 4585.46
$ws.Cells.Item(3, 29).Value = 691
$ws.Cells.Item(3, 30).Value = 29.95
$ws.Cells.Item(3, 31).Value = 7760
$ws.Cells.Item(3, 32).Value = 2.67
$ws.Cells.Item(3, 33).Value = 300
$ws.Cells.Item(3, 34).Value = 1.45
$ws.Cells.Item(3, 35).Value = 38.21
$ws.Cells.Item(3, 36).Value = 115041225

# Row 4
$ws.Cells.Item(4, 4).Value = 32326
$ws.Cells.Item(4, 5).Value = 1495
$ws.Cells.Item(4, 6).Value = 1495
$ws.Cells.Item(4, 7).Value = 1278
$ws.Cells.Item(4, 8).Value = 906
$ws.Cells.Item(4, 9).Value = 883
$ws.Cells.Item(4, 10).Value = 24
$ws.Cells.Item(4, 11).Value = 21514
$ws.Cells.Item(4, 12).Value = 13751
$ws.Cells.Item(4, 13).Value = 7763
$ws.Cells.Item(4, 14).Value = 7630
$ws.Cells.Item(4, 15).Value = 133
$ws.Cells.Item(4, 16).Value = 230
$ws.Cells.Item(4, 17).Value = 998
$ws.Cells.Item(4, 18).Value = -330
$ws.Cells.Item(4, 19).Value = -461
$ws.Cells.Item(4, 20).Value = 238
$ws.Cells.Item(4, 21).Value = 760
$ws.Cells.Item(4, 22).Value = 195
$ws.Cells.Item(4, 23).Value = 4.63
$ws.Cells.Item(4, 24).Value = 2.8
$ws.Cells.Item(4, 25).Value = 11.4
$ws.Cells.Item(4, 26).Value = 4.51
$ws.Cells.Item(4, 27).Value = 177.15
$ws.Cells.Item(4, 28).Value = 4868.37
$ws.Cells.Item(4, 29).Value = 767
$ws.Cells.Item(4, 30).Value = 20.53
$ws.Cells.Item(4, 31).Value = 7534
$ws.Cells.Item(4, 32).Value = 2.09
$ws.Cells.Item(4, 33).Value = 300
$ws.Cells.Item(4, 34).Value = 1.9
$ws.Cells.Item(4, 35).Value = 34.42
$ws.Cells.Item(4, 36).Value = 115041225

# Row 5
$ws.Cells.Item(5, 4).Value = 33750
$ws.Cells.Item(5, 5).Value = 1565
$ws.Cells.Item(5, 6).Value = 1565
$ws.Cells.Item(5, 7).Value = 1653
$ws.Cells.Item(5, 8).Value = 1284
$ws.Cells.Item(5, 9).Value = 1272
$ws.Cells.Item(5, 10).Value = 12
$ws.Cells.Item(5, 11).Value = 22383
$ws.Cells.Item(5, 12).Value = 13823
$ws.Cells.Item(5, 13).Value = 8560
$ws.Cells.Item(5, 14).Value = 8461
$ws.Cells.Item(5, 15).Value = 99
$ws.Cells.Item(5, 16).Value = 230
$ws.Cells.Item(5, 17).Value = 391
$ws.Cells.Item(5, 18).Value = -146
$ws.Cells.Item(5, 19).Value = -527
$ws.Cells.Item(5, 20).Value = 207
$ws.Cells.Item(5, 21).Value = 183
$ws.Cells.Item(5, 22).Value = 133
$ws.Cells.Item(5, 23).Value = 4.64
$ws.Cells.Item(5, 24).Value = 3.8
$ws.Cells.Item(5, 25).Value = 15.8
$ws.Cells.Item(5, 26).Value = 5.85
$ws.Cells.Item(5, 27).Value = 161.49
$ws.Cells.Item(5, 28).Value = 5290.76
$ws.Cells.Item(5, 29).Value = 1105
$ws.Cells.Item(5, 30).Value = 19.18
$ws.Cells.Item(5, 31).Value = 8354
$ws.Cells.Item(5, 32).Value = 2.54
$ws.Cells.Item(5, 33).Value = 760
$ws.Cells.Item(5, 34).Value = 3.58
$ws.Cells.Item(5, 35).Value = 60.53
$ws.Cells.Item(5, 36).Value = 115041225

# Row 6
$ws.Cells.Item(6, 4).Value = 34779
$ws.Cells.Item(6, 5).Value = 1811
$ws.Cells.Item(6, 6).Value = 1811
$ws.Cells.Item(6, 7).Value = 1905
$ws.Cells.Item(6, 8).Value = 1319
$ws.Cells.Item(6, 9).Value = 1297
$ws.Cells.Item(6, 11).Value = 21721
$ws.Cells.Item(6, 12).Value = 12793
$ws.Cells.Item(6, 13).Value = 8928
$ws.Cells.Item(6, 14).Value = 8826
$ws.Cells.Item(6, 16).Value = 230
$ws.Cells.Item(6, 17).Value = 2240
$ws.Cells.Item(6, 18).Value = -875
$ws.Cells.Item(6, 19).Value = -964
$ws.Cells.Item(6, 20).Value = 203
$ws.Cells.Item(6, 21).Value = 2037
$ws.Cells.Item(6, 22).Value = 176
$ws.Cells.Item(6, 23).Value = 5.21
$ws.Cells.Item(6, 24).Value = 3.79
$ws.Cells.Item(6, 25).Value = 15.01
$ws.Cells.Item(6, 26).Value = 5.98
$ws.Cells.Item(6, 27).Value = 143.29
$ws.Cells.Item(6, 28).Value = 5505.72
$ws.Cells.Item(6, 29).Value = 1128
$ws.Cells.Item(6, 30).Value = 19.96
$ws.Cells.Item(6, 31).Value = 8715
$ws.Cells.Item(6, 32).Value = 2.58
$ws.Cells.Item(6, 33).Value = 770
$ws.Cells.Item(6, 34).Value = 3.42
$ws.Cells.Item(6, 35).Value = 60.12
$ws.Cells.Item(6, 36).Value = 115041225

# Row 7
$ws.Cells.Item(7, 4).Value = 34838
$ws.Cells.Item(7, 5).Value = 2099
$ws.Cells.Item(7, 7).Value = 2153
$ws.Cells.Item(7, 8).Value = 1446
$ws.Cells.Item(7, 9).Value = 1430
$ws.Cells.Item(7, 11).Value = 23337
$ws.Cells.Item(7, 12).Value = 13738
$ws.Cells.Item(7, 13).Value = 9599
$ws.Cells.Item(7, 14).Value = 9490
$ws.Cells.Item(7, 16).Value = 230
$ws.Cells.Item(7, 17).Value = 1758
$ws.Cells.Item(7, 18).Value = -280
$ws.Cells.Item(7, 19).Value = -874
$ws.Cells.Item(7, 20).Value = 313
$ws.Cells.Item(7, 21).Value = 1202
$ws.Cells.Item(7, 23).Value = 6.03
$ws.Cells.Item(7, 24).Value = 4.15
$ws.Cells.Item(7, 25).Value = 15.61
$ws.Cells.Item(7, 26).Value = 6.42
$ws.Cells.Item(7, 27).Value = 143.12
$ws.Cells.Item(7, 29).Value = 1243
$ws.Cells.Item(7, 30).Value = 18.15
$ws.Cells.Item(7, 31).Value = 9370
$ws.Cells.Item(7, 32).Value = 2.41
$ws.Cells.Item(7, 33).Value = 809
$ws.Cells.Item(7, 34).Value = 3.59
$ws.Cells.Item(7, 35).Value = 65.06

# Row 8
$ws.Cells.Item(8, 4).Value = 36625
$ws.Cells.Item(8, 5).Value = 2261
$ws.Cells.Item(8, 7).Value = 2321
$ws.Cells.Item(8, 8).Value = 1607
$ws.Cells.Item(8, 9).Value = 1583
$ws.Cells.Item(8, 11).Value = 24376
$ws.Cells.Item(8, 12).Value = 14064
$ws.Cells.Item(8, 13).Value = 10312
$ws.Cells.Item(8, 14).Value = 10197
$ws.Cells.Item(8, 16).Value = 230
$ws.Cells.Item(8, 17).Value = 1870
$ws.Cells.Item(8, 18).Value = -394
$ws.Cells.Item(8, 19).Value = -893
$ws.Cells.Item(8, 20).Value = 226
$ws.Cells.Item(8, 21).Value = 1621
$ws.Cells.Item(8, 23).Value = 6.17
$ws.Cells.Item(8, 24).Value = 4.39
$ws.Cells.Item(8, 25).Value = 16.08
$ws.Cells.Item(8, 26).Value = 6.73
$ws.Cells.Item(8, 27).Value = 136.39
$ws.Cells.Item(8, 29).Value = 1376
$ws.Cells.Item(8, 30).Value = 15.37
$ws.Cells.Item(8, 31).Value = 10068
$ws.Cells.Item(8, 32).Value = 2.1
$ws.Cells.Item(8, 33).Value = 872
$ws.Cells.Item(8, 34).Value = 4.12
$ws.Cells.Item(8, 35).Value = 63.34

# Row 9
$ws.Cells.Item(9, 4).Value = 38753
$ws.Cells.Item(9, 5).Value = 2454
$ws.Cells.Item(9, 7).Value = 2503
$ws.Cells.Item(9, 8).Value = 1736
$ws.Cells.Item(9, 9).Value = 1711
$ws.Cells.Item(9, 11).Value = 25525
$ws.Cells.Item(9, 12).Value = 14405
$ws.Cells.Item(9, 13).Value = 11120
$ws.Cells.Item(9, 14).Value = 10977
$ws.Cells.Item(9, 16).Value = 230
$ws.Cells.Item(9, 17).Value = 1952
$ws.Cells.Item(9, 18).Value = -363
$ws.Cells.Item(9, 19).Value = -956
$ws.Cells.Item(9, 20).Value = 229
$ws.Cells.Item(9, 21).Value = 1546
$ws.Cells.Item(9, 23).Value = 6.33
$ws.Cells.Item(9, 24).Value = 4.48
$ws.Cells.Item(9, 25).Value = 16.16
$ws.Cells.Item(9, 26).Value = 6.96
$ws.Cells.Item(9, 27).Value = 129.55
$ws.Cells.Item(9, 29).Value = 1487
$ws.Cells.Item(9, 30).Value = 14.22
$ws.Cells.Item(9, 31).Value = 10838
$ws.Cells.Item(9, 32).Value = 1.95
$ws.Cells.Item(9, 33).Value = 921
$ws.Cells.Item(9, 34).Value = 4.36
$ws.Cells.Item(9, 35).Value = 61.95
